# Fruta / hortaliza, semanal
# A new daily-price record was inserted as row 314 in the "Naranja" sheet.
# All existing rows from 314 through 362 shift down by one (to 315-363),
# and the new row 314 is populated with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 314..362 down to 315..363, leaving a blank row 314.
$ws.Rows(314).Insert()

# Populate the newly inserted row 314 with the new record.
$ws.Cells.Item(314, 1).Value  = 5
$ws.Cells.Item(314, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(314, 3).Value  = "Maule"
$ws.Cells.Item(314, 4).Value  = 44522
$ws.Cells.Item(314, 5).Value  = 7
$ws.Cells.Item(314, 6).Value  = "Fruta"
$ws.Cells.Item(314, 7).Value  = 100102
$ws.Cells.Item(314, 8).Value  = "Cítricos"
$ws.Cells.Item(314, 9).Value  = 100102005
$ws.Cells.Item(314, 10).Value = "Naranja"
$ws.Cells.Item(314, 11).Value = "Lane Late"
$ws.Cells.Item(314, 12).Value = "Primera"
$ws.Cells.Item(314, 13).Value = 300
$ws.Cells.Item(314, 14).Value = 8000
$ws.Cells.Item(314, 15).Value = 8000
$ws.Cells.Item(314, 16).Value = 8000
$ws.Cells.Item(314, 17).Value = '$/bandeja 15 kilos granel'
$ws.Cells.Item(314, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(314, 19).Value = 533
$ws.Cells.Item(314, 20).Value = 15
